$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.540.60'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '1.920.74'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("E4").Value = '  +0.71%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.10'
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4818'
$ws.Range("E7").Value = '  -0.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4057'
$ws.Range("E8").Value = '  -0.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08235'
$ws.Range("E9").Value = '  +1.08%  '
$ws.Range("E10").Value = '  -0.15%  '
$ws.Range("E11").Value = '  -0.11%  '
$ws.Range("D12").Value = '1.931.05'
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("E13").Value = '  +0.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.256'
$ws.Range("E14").Value = '  +2.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.63'
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06846'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("E17").Value = '  +0.72%  '
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.55'
$ws.Range("E19").Value = '  -0.93%  '
$ws.Range("E20").Value = '  +0.74%  '
$ws.Range("D21").Value = '29.555.95'
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.677'
$ws.Range("E22").Value = '  +1.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.86'
$ws.Range("E23").Value = '  +0.57%  '
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("D25").Value = '2.184.39'
$ws.Range("E25").Value = '  +1.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.538'
$ws.Range("E26").Value = '  +4.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '155.99'
$ws.Range("E27").Value = '  +0.87%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.101'
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.68'
$ws.Range("E30").Value = '  +0.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.021'
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09641'
$ws.Range("E32").Value = '  +0.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.623'
$ws.Range("E33").Value = '  +1.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.559'
$ws.Range("E34").Value = '  +0.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.376'
$ws.Range("E35").Value = '  -1.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06308'
$ws.Range("E36").Value = '  +3.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02288'
$ws.Range("E37").Value = '  +0.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.186'
$ws.Range("E38").Value = '  +1.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5946'
$ws.Range("E39").Value = '  -0.11%  '
$ws.Range("E40").Value = '  +0.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.918'
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1849'
$ws.Range("E42").Value = '  -0.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.467'
$ws.Range("E43").Value = '  -0.52%  '
$ws.Range("E44").Value = '  -0.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.37'
$ws.Range("E45").Value = '  -0.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07470'
$ws.Range("E46").Value = '  -3.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5568'
$ws.Range("E47").Value = '  -0.11%  '
$ws.Range("E48").Value = '  -0.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '118.53'
$ws.Range("E49").Value = '  +3.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.431'
$ws.Range("E50").Value = '  +3.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.24'
$ws.Range("E51").Value = '  -0.68%  '
